$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (no value-type change) for numeric-looking Price values
# so Excel does not silently coerce them to Double and lose exact formatting
# (trailing zeros, grouping dots, etc.) - matches the original inlineStr text cells.

# Row 2
$ws.Range("D2").Value = '30.326.90'
$ws.Range("E2").Value = '  -3.13%  '

# Row 3
$ws.Range("D3").Value = '1.934.81'
$ws.Range("E3").Value = '  -3.56%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.45'
$ws.Range("E5").Value = '  -4.04%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7274'
$ws.Range("E6").Value = '  -7.49%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9989'
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3347'
$ws.Range("E8").Value = '  -7.16%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '28.68'
$ws.Range("E9").Value = '  +0.76%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07381'
$ws.Range("E10").Value = '  +4.28%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8159'
$ws.Range("E11").Value = '  -5.33%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08107'
$ws.Range("E12").Value = '  -1.03%  '

# Row 13
$ws.Range("D13").Value = '1.934.78'
$ws.Range("E13").Value = '  -3.59%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.484'
$ws.Range("E14").Value = '  -2.74%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '95.02'
$ws.Range("E15").Value = '  -6.40%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.84'
$ws.Range("E16").Value = '  -3.14%  '

# Row 17
$ws.Range("D17").Value = '30.327.93'
$ws.Range("E17").Value = '  -3.16%  '

# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008324'
$ws.Range("E18").Value = '  +3.80%  '

# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '254.78'
$ws.Range("E19").Value = '  -7.63%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.859'
$ws.Range("E20").Value = '  -1.70%  '

# Row 21
$ws.Range("D21").Value = '2.189.09'
$ws.Range("E21").Value = '  -3.51%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9988'
$ws.Range("E22").Value = '  -0.11%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9985'
$ws.Range("E23").Value = '  -0.21%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.942'
$ws.Range("E24").Value = '  -3.50%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.814'
$ws.Range("E25").Value = '  -3.13%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.19'
$ws.Range("E26").Value = '  -3.43%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.429'
$ws.Range("E27").Value = '  +1.26%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.29'
$ws.Range("E28").Value = '  -4.00%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1345'
$ws.Range("E29").Value = '  -10.92%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.560'
$ws.Range("E30").Value = '  -4.22%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.341'
$ws.Range("E31").Value = '  -1.33%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.441'
$ws.Range("E32").Value = '  -4.38%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.232'
$ws.Range("E33").Value = '  -4.97%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05196'
$ws.Range("E34").Value = '  -1.02%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.246'
$ws.Range("E35").Value = '  +1.69%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7512'
$ws.Range("E36").Value = '  -3.77%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.759'
$ws.Range("E37").Value = '  -1.92%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01997'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.834'
$ws.Range("E39").Value = '  -3.57%  '

# Row 40
$ws.Range("E40").Value = '  -1.14%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '79.28'
$ws.Range("E41").Value = '  -1.38%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4522'
$ws.Range("E42").Value = '  -5.17%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.020'
$ws.Range("E43").Value = '  -6.56%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9993'
$ws.Range("E44").Value = '  -0.11%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8392'
$ws.Range("E45").Value = '  -2.33%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.57'
$ws.Range("E46").Value = '  -5.07%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.799'
$ws.Range("E47").Value = '  -1.74%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.372'
$ws.Range("E48").Value = '  -5.84%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.81'
$ws.Range("E49").Value = '  -0.48%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.503'
$ws.Range("E50").Value = '  -0.02%  '

# Row 51
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4127'
$ws.Range("E51").Value = '  -5.67%  '
